# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and correct the Hedera/Filecoin row ordering (rows 40-41).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.632.81"
$ws.Range("E2").Value = "  +4.45%  "

$ws.Range("D3").Value = "3.495.60"
$ws.Range("E3").Value = "  +2.46%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "589.10"
$ws.Range("E5").Value = "  +3.42%  "

$ws.Range("D6").Value = "171.40"
$ws.Range("E6").Value = "  +7.85%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "3.494.02"
$ws.Range("E8").Value = "  +2.34%  "

$ws.Range("E9").Value = "  +5.04%  "

$ws.Range("E10").Value = "  +0.56%  "

$ws.Range("E11").Value = "  +4.84%  "

$ws.Range("D12").Value = "0.436"
$ws.Range("E12").Value = "  +3.44%  "

$ws.Range("D13").Value = "4.098.50"
$ws.Range("E13").Value = "  +2.42%  "

$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15").Value = "28.22"
$ws.Range("E15").Value = "  +4.41%  "

$ws.Range("D16").Value = "66.656.73"
$ws.Range("E16").Value = "  +4.38%  "

$ws.Range("E17").Value = "  +3.01%  "

$ws.Range("D18").Value = "3.474.37"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("D19").Value = "6.34"
$ws.Range("E19").Value = "  +4.61%  "

$ws.Range("E20").Value = "  +2.88%  "

$ws.Range("D21").Value = "388.30"
$ws.Range("E21").Value = "  +3.16%  "

$ws.Range("D22").Value = "7.94"
$ws.Range("E22").Value = "  +1.39%  "

$ws.Range("D23").Value = "73.09"
$ws.Range("E23").Value = "  +2.28%  "

$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("D25").Value = "0.532"
$ws.Range("E25").Value = "  +3.32%  "

$ws.Range("D26").Value = "0.0000122"
$ws.Range("E26").Value = "  +5.55%  "

$ws.Range("D27").Value = "10.39"
$ws.Range("E27").Value = "  +8.59%  "

$ws.Range("E28").Value = "  +2.87%  "

$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  +6.32%  "

$ws.Range("D31").Value = "1.48"
$ws.Range("E31").Value = "  +6.38%  "

$ws.Range("E32").Value = "  +2.46%  "

$ws.Range("E33").Value = "  +2.96%  "

$ws.Range("D34").Value = "7.43"
$ws.Range("E34").Value = "  +5.72%  "

$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("E36").Value = "  +7.01%  "

$ws.Range("D37").Value = "162.91"
$ws.Range("E37").Value = "  +2.12%  "

$ws.Range("D38").Value = "0.880"
$ws.Range("E38").Value = "  +5.77%  "

$ws.Range("E39").Value = "  +5.47%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "4.67"
$ws.Range("E40").Value = "  +5.66%  "

$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "0.0745"
$ws.Range("E41").Value = "  +2.13%  "

$ws.Range("D42").Value = "26.31"
$ws.Range("E42").Value = "  +2.31%  "

$ws.Range("D43").Value = "6.64"
$ws.Range("E43").Value = "  +2.42%  "

$ws.Range("D44").Value = "2.807.98"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("D45").Value = "26.65"
$ws.Range("E45").Value = "  +2.76%  "

$ws.Range("D46").Value = "42.95"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").Value = "356.51"
$ws.Range("E47").Value = "  +6.11%  "

$ws.Range("D48").Value = "2.52"
$ws.Range("E48").Value = "  +6.82%  "

$ws.Range("D49").Value = "0.0311"
$ws.Range("E49").Value = "  +2.75%  "

$ws.Range("D50").Value = "1.09"
$ws.Range("E50").Value = "  +4.48%  "

$ws.Range("D51").Value = "33.83"
$ws.Range("E51").Value = "  +14.07%  "
